# Streamline the wellplate import template so it can also be used for
# research plan table imports: the "sample_ID" column becomes a generic
# "Sample" column, and the pre-filled placeholder sample IDs (10001..10096)
# are cleared so the sheet ships blank and ready for either use case.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Wellplate import template")

# Rename header "sample_ID" -> "Sample"
$ws.Range("B1").Value = "Sample"

# Clear the pre-filled numeric sample IDs in column B (data rows 2-97)
$ws.Range("B2:B97").ClearContents()

# Move the active selection to C11 (matches the author's last worked cell)
$ws.Range("C11").Select() | Out-Null
